$d = $word.ActiveDocument
$newText = "Συμμετέχετε σε μία παγκόσμια καμπάνια για να παρατηρήσετε και να καταγράψετε τη φωτεινότητα των πιο αμυδρά ορατών άστρων σαν μέσο για την μέτρηση της Φωτορρύπανσης σε μία δεδομένη περιοχή. Με τον εντοπισμό και την παρατήρηση του  Αστερισμός μπότες στον νυχτερινό ουρανό καθώς και με την σύγκριση των ανωτέρω με τα διαγράμματα για τα μεγέθη των άστρων,  άνθρωποι από όλον τον κόσμο θα μάθουν πώς τα φώτα στην κοινότητά τους συμβάλλουν στην Φωτορρύπανση. Με την κατάθεση των πορισμάτων τους στην ιστοσελίδα θα δημιουργηθεί ένα αρχείο σχετικά με το τι μπορεί να δει κανείς στον νυχτερινό ουρανό."
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Συμμετέχετε*") {
        $r = $p.Range
        $r.MoveEnd(1, -1)  # exclude paragraph mark
        $r.Delete() | Out-Null
        $r.InsertAfter($newText) | Out-Null
        break
    }
}
